$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD2").Value = 151

$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.4
$ws.Range("J3").Value = 2.38
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 5.5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("W3").Value = 6
$ws.Range("AA3").Value = 15
$ws.Range("AC3").Value = 8
$ws.Range("AF3").Value = 67
$ws.Range("AJ3").Value = 51
$ws.Range("AK3").Value = 41
$ws.Range("AM3").Value = 451
$ws.Range("AO3").Value = 9
$ws.Range("AR3").Value = 51
$ws.Range("AS3").Value = 201
$ws.Range("AU3").Value = 9
$ws.Range("AV3").Value = 67
$ws.Range("AX3").Value = 29
$ws.Range("AZ3").Value = 101
$ws.Range("BB3").Value = 351

$ws.Range("M4").Value = 1.14
$ws.Range("O4").Value = 1.67

$ws.Range("M5").Value = 1.13
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.2

$ws.Range("M6").Value = 1.11
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 2.38

$ws.Range("M7").Value = 1.14
$ws.Range("O7").Value = 1.62

$ws.Range("M8").Value = 1.04
$ws.Range("O8").Value = 1.25
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03

$ws.Range("R10").Value = 1.65

$ws.Range("R11").Value = 1.65
